$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last three rows (old rows 8, 9, 10 -> classes 12, 13, 14)
$ws.Rows("8:10").Delete()

# New class labels (A2:A7)
$classes = @("[54, 61)", "[61, 68)", "[68, 75)", "[75, 82)", "[82, 89)", "[89, 96)")
for ($i = 0; $i -lt $classes.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $classes[$i]
}

# New numeric data for B2:F7 (Marca de clase, ni, hi, Ni, Hi)
$data = @(
    @(57.5, 4, 0.133, 4,  0.133),
    @(64.5, 5, 0.167, 9,  0.3),
    @(71.5, 2, 0.067, 11, 0.367),
    @(78.5, 8, 0.267, 19, 0.633),
    @(85.5, 8, 0.267, 27, 0.9),
    @(92.5, 2, 0.067, 29, 1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
}

Write-Output "done"
